# Refresh the cryptos price list (data for GitHub Actions scheduled update).
# Column D ("Price") and column E ("Volume(1h)") are plain text cells in the
# source workbook (t="inlineStr"), e.g. "64.143.84" or "  -0.23%  ", so we
# must keep them as text rather than let Excel reinterpret numeric-looking
# values (such as "550.89") as real numbers. A leading apostrophe forces
# Excel's COM layer to store the value as text, matching the original
# inlineStr/shared-string representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (column D) values ---
$ws.Range("D2").Value = "64.140.46"
$ws.Range("D3").Value = "3.323.57"
$ws.Range("D5").Value = "'550.89"
$ws.Range("D6").Value = "'172.83"
$ws.Range("D7").Value = "'0.625"
$ws.Range("D9").Value = "3.313.09"
$ws.Range("D12").Value = "'53.14"
$ws.Range("D13").Value = "'0.0000278"
$ws.Range("D15").Value = "3.854.35"
$ws.Range("D18").Value = "3.334.58"
$ws.Range("D19").Value = "63.998.44"
$ws.Range("D20").Value = "'11.68"
$ws.Range("D22").Value = "'447.98"
$ws.Range("D23").Value = "'4.97"
$ws.Range("D27").Value = "'2.86"
$ws.Range("D29").Value = "'8.55"
$ws.Range("D30").Value = "'30.74"
$ws.Range("D32").Value = "'62.57"
$ws.Range("D33").Value = "'11.33"
$ws.Range("D34").Value = "'570.03"
$ws.Range("D38").Value = "'3.54"
$ws.Range("D39").Value = "'35.11"
$ws.Range("D40").Value = "'0.365"
$ws.Range("D41").Value = "0.0₃0725"
$ws.Range("D42").Value = "3.053.09"
$ws.Range("D49").Value = "'142.02"
$ws.Range("D50").Value = "'2.51"
$ws.Range("D51").Value = "'8.16"

# --- Update Volume(1h) (column E) values ---
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E7").Value = "  +1.57%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("E10").Value = "  +6.07%  "
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("E16").Value = "  +2.07%  "
$ws.Range("E17").Value = "  -1.94%  "
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  +4.65%  "
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("E30").Value = "  +3.16%  "
$ws.Range("E31").Value = "  -2.97%  "
$ws.Range("E32").Value = "  +7.35%  "
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("E41").Value = "  -4.68%  "
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("E44").Value = "  -3.71%  "
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("E49").Value = "  +5.16%  "
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("E51").Value = "  -1.21%  "

# --- Row re-ranking swaps (B, C, D, E) ---
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'86.88"
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'13.86"
$ws.Range("E26").Value = "  +4.04%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.133"
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.43"
$ws.Range("E47").Value = "  -1.57%  "

Write-Host "Applied cryptos update."
